$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet
$ws.Name = "SoccerPage"

# Update A2 and add new team names in A3:A10
$ws.Range("A2").Value = "REAL MADRID"
$ws.Range("A3").Value = "PARIS SAINT-GERMAIN"
$ws.Range("A4").Value = "MANCHESTER UNITED"
$ws.Range("A5").Value = "MANCHESTER CITY"
$ws.Range("A6").Value = "JUVENTUS"
$ws.Range("A7").Value = "BAYERN MUNICH"
$ws.Range("A8").Value = "BARCELONA"
$ws.Range("A9").Value = "LIVERPOOL"
$ws.Range("A10").Value = "BORUSSIA DORTMUND"

# Set column A width to fit content (bestFit)
$ws.Columns.Item(1).AutoFit() | Out-Null

# Select the next empty cell, like Excel would leave the cursor after manual entry
$ws.Range("A11").Select() | Out-Null
